$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value2 = 5299.9  # H32
$ws.Cells.Item(32, 10).Value2 = 5749.75  # J32
$ws.Cells.Item(32, 12).Value2 = 5749.75  # L32
$ws.Cells.Item(32, 14).Value2 = -6401.75  # N32
$ws.Cells.Item(40, 8).Value2 = 3321.111  # H40
$ws.Cells.Item(40, 9).Value2 = 3618.889  # I40
$ws.Cells.Item(40, 11).Value2 = 3618.889  # K40
$ws.Cells.Item(40, 13).Value2 = -3443.889  # M40
$ws.Cells.Item(62, 8).Value2 = 4395.5625  # H62
$ws.Cells.Item(62, 9).Value2 = 4050  # I62
$ws.Cells.Item(62, 10).Value2 = 5155.8  # J62
$ws.Cells.Item(62, 11).Value2 = 4050  # K62
$ws.Cells.Item(62, 12).Value2 = 5155.8  # L62
$ws.Cells.Item(62, 13).Value2 = -3426  # M62
$ws.Cells.Item(62, 14).Value2 = -6403.8  # N62
$ws.Cells.Item(65, 8).Value2 = 4395.5625  # H65
$ws.Cells.Item(65, 9).Value2 = 4050  # I65
$ws.Cells.Item(65, 10).Value2 = 5155.8  # J65
$ws.Cells.Item(65, 11).Value2 = 20250  # K65
$ws.Cells.Item(65, 12).Value2 = 25779  # L65
$ws.Cells.Item(65, 13).Value2 = -17130  # M65
$ws.Cells.Item(65, 14).Value2 = -32019  # N65
$ws.Cells.Item(70, 8).Value2 = 0  # H70
$ws.Cells.Item(70, 9).Value2 = 0  # I70
$ws.Cells.Item(70, 10).Value2 = 0  # J70
$ws.Cells.Item(70, 11).Value2 = 0  # K70
$ws.Cells.Item(70, 12).Value2 = 0  # L70
$ws.Cells.Item(73, 8).Value2 = 0  # H73
$ws.Cells.Item(73, 9).Value2 = 0  # I73
$ws.Cells.Item(73, 10).Value2 = 0  # J73
$ws.Cells.Item(73, 11).Value2 = 0  # K73
$ws.Cells.Item(73, 12).Value2 = 0  # L73
$ws.Cells.Item(74, 8).Value2 = 3710.889  # H74
$ws.Cells.Item(74, 10).Value2 = 0  # J74
$ws.Cells.Item(74, 12).Value2 = 0  # L74
$ws.Cells.Item(77, 8).Value2 = 3710.889  # H77
$ws.Cells.Item(77, 10).Value2 = 0  # J77
$ws.Cells.Item(77, 12).Value2 = 0  # L77
$ws.Cells.Item(137, 8).Value2 = 191208.89  # H137
$ws.Cells.Item(137, 9).Value2 = 257265.52  # I137
$ws.Cells.Item(137, 10).Value2 = 7194  # J137
$ws.Cells.Item(137, 11).Value2 = 771796.5599999999  # K137
$ws.Cells.Item(137, 12).Value2 = 21582  # L137
$ws.Cells.Item(137, 13).Value2 = -769246.5599999999  # M137
$ws.Cells.Item(137, 14).Value2 = -26682  # N137
$ws.Cells.Item(70, 13).ClearContents()  # M70
$ws.Cells.Item(70, 14).ClearContents()  # N70
$ws.Cells.Item(73, 13).ClearContents()  # M73
$ws.Cells.Item(73, 14).ClearContents()  # N73
$ws.Cells.Item(74, 14).ClearContents()  # N74
$ws.Cells.Item(77, 14).ClearContents()  # N77

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value2 = 5499846.5  # H32
$ws.Cells.Item(32, 9).Value2 = 6668726.5  # I32
$ws.Cells.Item(32, 11).Value2 = 6668726.5  # K32
$ws.Cells.Item(32, 13).Value2 = -6668439.5  # M32
$ws.Cells.Item(45, 8).Value2 = 71430590  # H45
$ws.Cells.Item(45, 9).Value2 = 83335020  # I45
$ws.Cells.Item(45, 11).Value2 = 83335020  # K45
$ws.Cells.Item(45, 13).Value2 = -83334643  # M45
$ws.Cells.Item(74, 8).Value2 = 6763562  # H74
$ws.Cells.Item(74, 9).Value2 = 10002040  # I74
$ws.Cells.Item(74, 10).Value2 = 16734.334  # J74
$ws.Cells.Item(74, 11).Value2 = 10002040  # K74
$ws.Cells.Item(74, 12).Value2 = 16734.334  # L74
$ws.Cells.Item(74, 13).Value2 = -10001166  # M74
$ws.Cells.Item(74, 14).Value2 = -18482.334  # N74
$ws.Cells.Item(77, 8).Value2 = 6763562  # H77
$ws.Cells.Item(77, 9).Value2 = 10002040  # I77
$ws.Cells.Item(77, 10).Value2 = 16734.334  # J77
$ws.Cells.Item(77, 11).Value2 = 50010200  # K77
$ws.Cells.Item(77, 12).Value2 = 83671.67  # L77
$ws.Cells.Item(77, 13).Value2 = -50005832  # M77
$ws.Cells.Item(77, 14).Value2 = -92407.67  # N77

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(36, 8).Value2 = 3180.8333  # H36
$ws.Cells.Item(36, 10).Value2 = 4977.5  # J36
$ws.Cells.Item(36, 12).Value2 = 4977.5  # L36
$ws.Cells.Item(36, 14).Value2 = -6045.5  # N36
$ws.Cells.Item(82, 8).Value2 = 30857.6  # H82
$ws.Cells.Item(82, 9).Value2 = 30857.6  # I82
$ws.Cells.Item(82, 11).Value2 = 30857.6  # K82
$ws.Cells.Item(82, 13).Value2 = -30474.6  # M82
$ws.Cells.Item(85, 8).Value2 = 30857.6  # H85
$ws.Cells.Item(85, 9).Value2 = 30857.6  # I85
$ws.Cells.Item(85, 11).Value2 = 30857.6  # K85
$ws.Cells.Item(85, 13).Value2 = -29531.6  # M85
$ws.Cells.Item(94, 8).Value2 = 1402.8518  # H94
$ws.Cells.Item(94, 9).Value2 = 518.4  # I94
$ws.Cells.Item(94, 11).Value2 = 518.4  # K94
$ws.Cells.Item(94, 13).Value2 = -67.39999999999998  # M94

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value2 = 1376.8572  # H16
$ws.Cells.Item(16, 9).Value2 = 1232  # I16
$ws.Cells.Item(16, 10).Value2 = 1485.5  # J16
$ws.Cells.Item(16, 11).Value2 = 1232  # K16
$ws.Cells.Item(16, 12).Value2 = 1485.5  # L16
$ws.Cells.Item(16, 13).Value2 = -945  # M16
$ws.Cells.Item(16, 14).Value2 = -2059.5  # N16
$ws.Cells.Item(105, 8).Value2 = 1977.091  # H105
$ws.Cells.Item(105, 9).Value2 = 1973.7  # I105
$ws.Cells.Item(105, 11).Value2 = 1973.7  # K105
$ws.Cells.Item(105, 13).Value2 = -226.7  # M105
$ws.Cells.Item(109, 8).Value2 = 0  # H109
$ws.Cells.Item(109, 10).Value2 = 0  # J109
$ws.Cells.Item(109, 12).Value2 = 0  # L109
$ws.Cells.Item(113, 8).Value2 = 1376.8572  # H113
$ws.Cells.Item(113, 9).Value2 = 1232  # I113
$ws.Cells.Item(113, 10).Value2 = 1485.5  # J113
$ws.Cells.Item(113, 11).Value2 = 1232  # K113
$ws.Cells.Item(113, 12).Value2 = 1485.5  # L113
$ws.Cells.Item(113, 13).Value2 = 938  # M113
$ws.Cells.Item(113, 14).Value2 = -5825.5  # N113
$ws.Cells.Item(109, 14).ClearContents()  # N109

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(55, 8).Value2 = 4379.8  # H55
$ws.Cells.Item(55, 9).Value2 = 3724.75  # I55
$ws.Cells.Item(55, 11).Value2 = 11174.25  # K55
$ws.Cells.Item(55, 13).Value2 = -10997.25  # M55
$ws.Cells.Item(68, 8).Value2 = 2866.2  # H68
$ws.Cells.Item(68, 9).Value2 = 1999.6666  # I68
$ws.Cells.Item(68, 10).Value2 = 4166  # J68
$ws.Cells.Item(68, 11).Value2 = 5998.9998  # K68
$ws.Cells.Item(68, 12).Value2 = 12498  # L68
$ws.Cells.Item(68, 13).Value2 = -5187.9998  # M68
$ws.Cells.Item(68, 14).Value2 = -14120  # N68
$ws.Cells.Item(71, 8).Value2 = 2866.2  # H71
$ws.Cells.Item(71, 9).Value2 = 1999.6666  # I71
$ws.Cells.Item(71, 10).Value2 = 4166  # J71
$ws.Cells.Item(71, 11).Value2 = 17996.9994  # K71
$ws.Cells.Item(71, 12).Value2 = 37494  # L71
$ws.Cells.Item(71, 13).Value2 = -13940.9994  # M71
$ws.Cells.Item(71, 14).Value2 = -45606  # N71
$ws.Cells.Item(113, 8).Value2 = 1451.6875  # H113
$ws.Cells.Item(113, 9).Value2 = 735.8  # I113
$ws.Cells.Item(113, 10).Value2 = 1777.091  # J113
$ws.Cells.Item(113, 11).Value2 = 2207.4  # K113
$ws.Cells.Item(113, 12).Value2 = 5331.272999999999  # L113
$ws.Cells.Item(113, 13).Value2 = -37.39999999999964  # M113
$ws.Cells.Item(113, 14).Value2 = -9671.272999999999  # N113

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value2 = 2138.8215  # H102
$ws.Cells.Item(102, 9).Value2 = 1291.9333  # I102
$ws.Cells.Item(102, 10).Value2 = 3116  # J102
$ws.Cells.Item(102, 11).Value2 = 1291.9333  # K102
$ws.Cells.Item(102, 12).Value2 = 3116  # L102
$ws.Cells.Item(102, 13).Value2 = 330.0667000000001  # M102
$ws.Cells.Item(102, 14).Value2 = -6360  # N102
$ws.Cells.Item(122, 8).Value2 = 2345  # H122
$ws.Cells.Item(122, 9).Value2 = 1952.6666  # I122
$ws.Cells.Item(122, 11).Value2 = 5857.9998  # K122
$ws.Cells.Item(122, 13).Value2 = -3407.9998  # M122
$ws.Cells.Item(132, 8).Value2 = 24392552  # H132
$ws.Cells.Item(132, 9).Value2 = 25643324  # I132
$ws.Cells.Item(132, 11).Value2 = 76929972  # K132
$ws.Cells.Item(132, 13).Value2 = -76927442  # M132

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(4, 8).Value2 = 0  # H4
$ws.Cells.Item(4, 10).Value2 = 0  # J4
$ws.Cells.Item(4, 12).Value2 = 0  # L4
$ws.Cells.Item(28, 8).Value2 = 0  # H28
$ws.Cells.Item(28, 10).Value2 = 0  # J28
$ws.Cells.Item(28, 12).Value2 = 0  # L28
$ws.Cells.Item(37, 8).Value2 = 0  # H37
$ws.Cells.Item(37, 10).Value2 = 0  # J37
$ws.Cells.Item(37, 12).Value2 = 0  # L37
$ws.Cells.Item(45, 8).Value2 = 0  # H45
$ws.Cells.Item(45, 10).Value2 = 0  # J45
$ws.Cells.Item(45, 12).Value2 = 0  # L45
$ws.Cells.Item(46, 8).Value2 = 2705.5  # H46
$ws.Cells.Item(46, 9).Value2 = 2089.111  # I46
$ws.Cells.Item(46, 10).Value2 = 3321.889  # J46
$ws.Cells.Item(46, 11).Value2 = 2089.111  # K46
$ws.Cells.Item(46, 12).Value2 = 3321.889  # L46
$ws.Cells.Item(46, 13).Value2 = -1901.111  # M46
$ws.Cells.Item(46, 14).Value2 = -3697.889  # N46
$ws.Cells.Item(122, 8).Value2 = 5360.1816  # H122
$ws.Cells.Item(122, 9).Value2 = 4899.2  # I122
$ws.Cells.Item(122, 10).Value2 = 5744.3335  # J122
$ws.Cells.Item(122, 11).Value2 = 14697.6  # K122
$ws.Cells.Item(122, 12).Value2 = 17233.0005  # L122
$ws.Cells.Item(122, 13).Value2 = -12247.6  # M122
$ws.Cells.Item(122, 14).Value2 = -22133.0005  # N122
$ws.Cells.Item(4, 14).ClearContents()  # N4
$ws.Cells.Item(28, 14).ClearContents()  # N28
$ws.Cells.Item(37, 14).ClearContents()  # N37
$ws.Cells.Item(45, 14).ClearContents()  # N45

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(12, 8).Value2 = 22000000  # H12
$ws.Cells.Item(12, 10).Value2 = 22000000  # J12
$ws.Cells.Item(12, 12).Value2 = 22000000  # L12
$ws.Cells.Item(12, 14).Value2 = -22000284  # N12
$ws.Cells.Item(107, 8).Value2 = 29413364  # H107
$ws.Cells.Item(107, 9).Value2 = 41668150  # I107
$ws.Cells.Item(107, 10).Value2 = 1877.6  # J107
$ws.Cells.Item(107, 11).Value2 = 125004450  # K107
$ws.Cells.Item(107, 12).Value2 = 5632.799999999999  # L107
$ws.Cells.Item(107, 13).Value2 = -125002530  # M107
$ws.Cells.Item(107, 14).Value2 = -9472.799999999999  # N107
$ws.Cells.Item(126, 8).Value2 = 1702.9286  # H126
$ws.Cells.Item(126, 9).Value2 = 1486.75  # I126
$ws.Cells.Item(126, 11).Value2 = 4460.25  # K126
$ws.Cells.Item(126, 13).Value2 = -1990.25  # M126
$ws.Cells.Item(127, 8).Value2 = 63996.668  # H127
$ws.Cells.Item(127, 10).Value2 = 63996.668  # J127
$ws.Cells.Item(127, 12).Value2 = 63996.668  # L127
$ws.Cells.Item(127, 14).Value2 = -73916.66800000001  # N127
